# Auto-generated edit script applying the Ixion_Profits.xlsx diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 6851248.5
$ws.Range("I40").Value = 1974.5254
$ws.Range("J40").Value = 35716050
$ws.Range("K40").Value = 1974.5254
$ws.Range("L40").Value = 35716050
$ws.Range("M40").Value = -1799.5254
$ws.Range("N40").Value = -35716400

# Row 51: A Bile Business | Shark Oil
$ws.Range("H51").Value = 1576.2941
$ws.Range("I51").Value = 1400.1428
$ws.Range("J51").Value = 1699.6
$ws.Range("K51").Value = 1400.1428
$ws.Range("L51").Value = 1699.6
$ws.Range("M51").Value = -916.1428000000001
$ws.Range("N51").Value = -2667.6

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 1327.7084
$ws.Range("I98").Value = 1263.6957
$ws.Range("J98").Value = 2800
$ws.Range("K98").Value = 1263.6957
$ws.Range("L98").Value = 2800
$ws.Range("M98").Value = 234.3043
$ws.Range("N98").Value = -5796

# Row 115: 5-bell Energy | Competent Craftsman's Syrup
$ws.Range("H115").Value = 935.625
$ws.Range("I115").Value = 655
$ws.Range("J115").Value = 2900
$ws.Range("K115").Value = 1965
$ws.Range("L115").Value = 8700
$ws.Range("M115").Value = -398
$ws.Range("N115").Value = -11834

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 1327.7084
$ws.Range("I122").Value = 1263.6957
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 3791.0871
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -1341.0871
$ws.Range("N122").Value = -13300

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 1694.3889
$ws.Range("I132").Value = 1248.129
$ws.Range("J132").Value = 4461.2
$ws.Range("K132").Value = 3744.387
$ws.Range("L132").Value = 13383.6
$ws.Range("M132").Value = -1214.387
$ws.Range("N132").Value = -18443.6

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 1350.1428
$ws.Range("I137").Value = 1400.1818
$ws.Range("J137").Value = 1166.6666
$ws.Range("K137").Value = 4200.5454
$ws.Range("L137").Value = 3499.9998
$ws.Range("M137").Value = -1650.5454
$ws.Range("N137").Value = -8599.9998

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 722.2727
$ws.Range("I2").Value = 537.44446
$ws.Range("J2").Value = 1554
$ws.Range("K2").Value = 537.44446
$ws.Range("L2").Value = 1554
$ws.Range("M2").Value = -424.44446
$ws.Range("N2").Value = -1780

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 288797.44
$ws.Range("I61").Value = 2359.8125
$ws.Range("J61").Value = 530008.0600000001
$ws.Range("K61").Value = 2359.8125
$ws.Range("L61").Value = 530008.0600000001
$ws.Range("M61").Value = -2147.8125
$ws.Range("N61").Value = -530432.0600000001

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1242.5178
$ws.Range("I74").Value = 881.5143
$ws.Range("J74").Value = 1844.1904
$ws.Range("K74").Value = 881.5143
$ws.Range("L74").Value = 1844.1904
$ws.Range("M74").Value = -7.514300000000048
$ws.Range("N74").Value = -3592.1904

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1242.5178
$ws.Range("I77").Value = 881.5143
$ws.Range("J77").Value = 1844.1904
$ws.Range("K77").Value = 4407.5715
$ws.Range("L77").Value = 9220.951999999999
$ws.Range("M77").Value = -39.57150000000001
$ws.Range("N77").Value = -17956.952

# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 722.2727
$ws.Range("I116").Value = 537.44446
$ws.Range("J116").Value = 1554
$ws.Range("K116").Value = 537.44446
$ws.Range("L116").Value = 1554
$ws.Range("M116").Value = 1756.55554
$ws.Range("N116").Value = -6142

# Row 123: The Armoire Is Open | High Durium Armguards of Maiming
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 288797.44
$ws.Range("I136").Value = 2359.8125
$ws.Range("J136").Value = 530008.0600000001
$ws.Range("K136").Value = 7079.4375
$ws.Range("L136").Value = 1590024.18
$ws.Range("M136").Value = -4529.4375
$ws.Range("N136").Value = -1595124.18

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 722.2727
$ws.Range("I3").Value = 537.44446
$ws.Range("J3").Value = 1554
$ws.Range("K3").Value = 537.44446
$ws.Range("L3").Value = 1554
$ws.Range("M3").Value = -423.44446
$ws.Range("N3").Value = -1782

# Row 140: Ceremonial Teeth | Ra'Kaznar Twinfangs
$ws.Range("H140").Value = 39061.395
$ws.Range("J140").Value = 39061.395
$ws.Range("L140").Value = 39061.395
$ws.Range("N140").Value = -49421.395

$ws = $wb.Worksheets.Item("CRP")
# Row 13: Compulsory Conjury | Maple Cane
$ws.Range("H13").Value = 85005
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 85005
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 85005
$ws.Range("N13").Value = -85283
$ws.Range("M13").ClearContents()

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 3105.3433
$ws.Range("I31").Value = 2017.9688
$ws.Range("J31").Value = 4099.514
$ws.Range("K31").Value = 2017.9688
$ws.Range("L31").Value = 4099.514
$ws.Range("M31").Value = -1722.9688
$ws.Range("N31").Value = -4689.514

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 3105.3433
$ws.Range("I34").Value = 2017.9688
$ws.Range("J34").Value = 4099.514
$ws.Range("K34").Value = 2017.9688
$ws.Range("L34").Value = 4099.514
$ws.Range("M34").Value = -1815.9688
$ws.Range("N34").Value = -4503.514

# Row 39: An Expected Tourney | Ash Cavalry Bow
$ws.Range("H39").Value = 5025.5
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Row 49: Bend It Like Durendaire | Ash Cavalry Bow
$ws.Range("H49").Value = 5025.5
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 7823279
$ws.Range("I99").Value = 20250
$ws.Range("K99").Value = 20250
$ws.Range("M99").Value = -18752

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 1699.12
$ws.Range("I105").Value = 1674.9375
$ws.Range("J105").Value = 1742.1111
$ws.Range("K105").Value = 1674.9375
$ws.Range("L105").Value = 1742.1111
$ws.Range("M105").Value = 72.0625
$ws.Range("N105").Value = -5236.1111

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 7823279
$ws.Range("I126").Value = 20250
$ws.Range("K126").Value = 60750
$ws.Range("M126").Value = -58280

$ws = $wb.Worksheets.Item("CUL")
# Row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 1316318.6
$ws.Range("I113").Value = 2083866.1
$ws.Range("J113").Value = 523.0714
$ws.Range("K113").Value = 6251598.300000001
$ws.Range("L113").Value = 1569.2142
$ws.Range("M113").Value = -6249428.300000001
$ws.Range("N113").Value = -5909.2142

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 2241.647
$ws.Range("I132").Value = 2574.1667
$ws.Range("J132").Value = 2060.2727
$ws.Range("K132").Value = 23167.5003
$ws.Range("L132").Value = 18542.4543
$ws.Range("M132").Value = -20637.5003
$ws.Range("N132").Value = -23602.4543

# Row 139: Najoothie | Wild Banana Blend
$ws.Range("H139").Value = 4900.278
$ws.Range("I139").Value = 7275.625
$ws.Range("K139").Value = 21826.875
$ws.Range("M139").Value = -16686.875

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 425287.4
$ws.Range("I102").Value = 998030.25
$ws.Range("K102").Value = 998030.25
$ws.Range("M102").Value = -996408.25

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 2217.2646
$ws.Range("I122").Value = 1606.6923
$ws.Range("J122").Value = 2595.238
$ws.Range("K122").Value = 4820.0769
$ws.Range("L122").Value = 7785.714
$ws.Range("M122").Value = -2370.0769
$ws.Range("N122").Value = -12685.714

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 3235.5625
$ws.Range("I132").Value = 3919.8667
$ws.Range("J132").Value = 2631.7646
$ws.Range("K132").Value = 11759.6001
$ws.Range("L132").Value = 7895.293799999999
$ws.Range("M132").Value = -9229.6001
$ws.Range("N132").Value = -12955.2938

# Row 134: Guaranteed Gem | Ihuykanite
$ws.Range("H134").Value = 30500
$ws.Range("J134").Value = 30500
$ws.Range("L134").Value = 91500
$ws.Range("N134").Value = -96570

$ws = $wb.Worksheets.Item("LTW")
# Row 135: Dreams of Ja | Crocodileskin Leg Wraps of Scouting
$ws.Range("H135").Value = 42400
$ws.Range("J135").Value = 42400
$ws.Range("L135").Value = 42400
$ws.Range("N135").Value = -52540

# Row 141: Just Generally Freezing | Gargantuaskin Trousers of Striking
$ws.Range("H141").Value = 56123
$ws.Range("J141").Value = 56123
$ws.Range("L141").Value = 56123
$ws.Range("N141").Value = -66483

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1516.9778
$ws.Range("I132").Value = 1114.5625
$ws.Range("J132").Value = 2507.5386
$ws.Range("K132").Value = 3343.6875
$ws.Range("L132").Value = 7522.6158
$ws.Range("M132").Value = -813.6875
$ws.Range("N132").Value = -12582.6158
